$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: Row, Fecha(D), Calidad(I), Volumen(J), PrecioMin(K), PrecioMax(L), PrecioProm(M), Origen(O), PrecioKg(P)
$rows = @(
    @{Row=2; D=44855; I='Primera'; J=7900; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=3; D=44166; I='Primera'; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=4; D=44189; I='Primera'; J=16000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=5; D=44600; I='Primera'; J=1300; K=3500; L=4000; M=3808; O='Región Metropolitana'; P=38}
    @{Row=6; D=44882; I='Primera'; J=7900; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=7; D=44187; I='Primera'; J=12000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=8; D=44215; I='Primera'; J=16000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=9; D=44230; I='Primera'; J=16000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=10; D=44602; I='Primera'; J=12000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=11; D=44602; I='Segunda'; J=6000; K=2500; L=2500; M=2500; O='Provincia de Chacabuco'; P=25}
    @{Row=12; D=44229; I='Primera'; J=16000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=13; D=44159; I='Primera'; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=14; D=44902; I='Primera'; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=15; D=44214; I='Primera'; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=16; D=44875; I='Primera'; J=7900; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=17; D=44209; I='Primera'; J=7000; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=18; D=44883; I='Primera'; J=9700; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=19; D=44168; I='Primera'; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=20; D=44161; I='Primera'; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=21; D=44860; I='Primera'; J=7900; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=22; D=44210; I='Primera'; J=8800; K=2500; L=3000; M=2750; O='Provincia de Chacabuco'; P=28}
    @{Row=23; D=44859; I='Primera'; J=7900; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=24; D=44186; I='Primera'; J=10000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=25; D=44245; I='Primera'; J=9000; K=3000; L=3000; M=3000; O='Región Metropolitana'; P=30}
    @{Row=26; D=44245; I='Segunda'; J=5000; K=2500; L=2500; M=2500; O='Región Metropolitana'; P=25}
    @{Row=27; D=44876; I='Primera'; J=7900; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=28; D=44873; I='Primera'; J=7900; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=29; D=44845; I='Primera'; J=7900; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=30; D=44874; I='Primera'; J=7900; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=31; D=44880; I='Primera'; J=7900; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=32; D=44901; I='Primera'; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=33; D=44232; I='Primera'; J=16000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=34; D=44160; I='Primera'; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=35; D=44881; I='Primera'; J=7900; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=36; D=44846; I='Primera'; J=7900; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=37; D=44181; I='Primera'; J=12000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=38; D=44162; I='Primera'; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=39; D=44167; I='Primera'; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=40; D=44204; I='Primera'; J=7000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=41; D=44231; I='Primera'; J=12000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=42; D=44188; I='Primera'; J=12000; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
    @{Row=43; D=44847; I='Primera'; J=7900; K=3000; L=3000; M=3000; O='Provincia de Chacabuco'; P=30}
)

foreach ($rec in $rows) {
    $r = $rec.Row
    $ws.Range("D$r").Value = $rec.D
    $ws.Range("I$r").Value = $rec.I
    $ws.Range("J$r").Value = $rec.J
    $ws.Range("K$r").Value = $rec.K
    $ws.Range("L$r").Value = $rec.L
    $ws.Range("M$r").Value = $rec.M
    $ws.Range("O$r").Value = $rec.O
    $ws.Range("P$r").Value = $rec.P
}
